$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": status text for the two source files changed.
#  - 56da58e0-...md row (row 2): status text itself changed from
#    "Handed back: in sync with en-US" to "In Translation"
#  - e027af39-...md row (row 3): status moved on to "Ready for handoff"
#  - Latest HO Xliff Generate Date (col G) refreshed for both rows
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("G2").Value = "2016-10-19 17:40:52"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-19 17:40:52"

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---------------------------------------------------------------------------
# Sheet "zh-cn": status -> Ready for handoff, refreshed handoff datetime,
# and a new handback-is-stale error message.
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-10-19 17:40:41"
$wsZhCn.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d37ece8868101fdbda670653fc139ead4a0edcd/e2e/56da58e0-d409-4166-a50b-3236372abc9f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7ee148d4e8bd225a4575b136276e88f2550a306/e2e/56da58e0-d409-4166-a50b-3236372abc9f.md."

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-10-19 17:40:41"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d37ece8868101fdbda670653fc139ead4a0edcd/e2e/e027af39-0bd2-4a8f-880c-cfe01c963cba.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7ee148d4e8bd225a4575b136276e88f2550a306/e2e/e027af39-0bd2-4a8f-880c-cfe01c963cba.md."

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# Sheet "de-de": status -> Ready for handoff, and a new handback-is-stale
# error message (handoff datetime unchanged on this sheet).
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d37ece8868101fdbda670653fc139ead4a0edcd/e2e/56da58e0-d409-4166-a50b-3236372abc9f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7ee148d4e8bd225a4575b136276e88f2550a306/e2e/56da58e0-d409-4166-a50b-3236372abc9f.md."

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d37ece8868101fdbda670653fc139ead4a0edcd/e2e/e027af39-0bd2-4a8f-880c-cfe01c963cba.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7ee148d4e8bd225a4575b136276e88f2550a306/e2e/e027af39-0bd2-4a8f-880c-cfe01c963cba.md."

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
